$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 888, shifting existing data (rows 888:955) down to (892:959).
$ws.Range("A888:R891").Insert()

# Populate the 4 newly inserted rows with this week's data.

# Row 888
$ws.Cells.Item(888, 1).Value = 1
$ws.Cells.Item(888, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(888, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(888, 4).Value = 44783
$ws.Cells.Item(888, 5).Value = 15
$ws.Cells.Item(888, 6).Value = 100112033
$ws.Cells.Item(888, 7).Value = "Lechuga"
$ws.Cells.Item(888, 8).Value = "Escarola"
$ws.Cells.Item(888, 9).Value = "Primera"
$ws.Cells.Item(888, 10).Value = 120
$ws.Cells.Item(888, 11).Value = 5000
$ws.Cells.Item(888, 12).Value = 6000
$ws.Cells.Item(888, 13).Value = 5500
$ws.Cells.Item(888, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(888, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(888, 16).Value = 458
$ws.Cells.Item(888, 17).Value = 12
$ws.Cells.Item(888, 18).Value = "Hortaliza"

# Row 889
$ws.Cells.Item(889, 1).Value = 1
$ws.Cells.Item(889, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(889, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(889, 4).Value = 44783
$ws.Cells.Item(889, 5).Value = 15
$ws.Cells.Item(889, 6).Value = 100112033
$ws.Cells.Item(889, 7).Value = "Lechuga"
$ws.Cells.Item(889, 8).Value = "Escarola"
$ws.Cells.Item(889, 9).Value = "Segunda"
$ws.Cells.Item(889, 10).Value = 120
$ws.Cells.Item(889, 11).Value = 5000
$ws.Cells.Item(889, 12).Value = 6000
$ws.Cells.Item(889, 13).Value = 5500
$ws.Cells.Item(889, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(889, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(889, 16).Value = 306
$ws.Cells.Item(889, 17).Value = 18
$ws.Cells.Item(889, 18).Value = "Hortaliza"

# Row 890
$ws.Cells.Item(890, 1).Value = 1
$ws.Cells.Item(890, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(890, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(890, 4).Value = 44783
$ws.Cells.Item(890, 5).Value = 15
$ws.Cells.Item(890, 6).Value = 100112033
$ws.Cells.Item(890, 7).Value = "Lechuga"
$ws.Cells.Item(890, 8).Value = "Marina"
$ws.Cells.Item(890, 9).Value = "Primera"
$ws.Cells.Item(890, 10).Value = 130
$ws.Cells.Item(890, 11).Value = 4000
$ws.Cells.Item(890, 12).Value = 5000
$ws.Cells.Item(890, 13).Value = 4500
$ws.Cells.Item(890, 14).Value = "$/caja 12 unidades"
$ws.Cells.Item(890, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(890, 16).Value = 375
$ws.Cells.Item(890, 17).Value = 12
$ws.Cells.Item(890, 18).Value = "Hortaliza"

# Row 891
$ws.Cells.Item(891, 1).Value = 1
$ws.Cells.Item(891, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(891, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(891, 4).Value = 44783
$ws.Cells.Item(891, 5).Value = 15
$ws.Cells.Item(891, 6).Value = 100112033
$ws.Cells.Item(891, 7).Value = "Lechuga"
$ws.Cells.Item(891, 8).Value = "Marina"
$ws.Cells.Item(891, 9).Value = "Segunda"
$ws.Cells.Item(891, 10).Value = 140
$ws.Cells.Item(891, 11).Value = 4000
$ws.Cells.Item(891, 12).Value = 5000
$ws.Cells.Item(891, 13).Value = 4500
$ws.Cells.Item(891, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(891, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(891, 16).Value = 250
$ws.Cells.Item(891, 17).Value = 18
$ws.Cells.Item(891, 18).Value = "Hortaliza"

# Ensure the date column (D) in the new rows uses the same date format as the rest of the column.
$ws.Range("D888:D891").NumberFormat = $ws.Range("D892").NumberFormat
